$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.83"
$ws.Range("E2").Value = "'3.95%"
$ws.Range("G2").Value = "'16"

$ws.Range("D3").Value = "'35.55"
$ws.Range("E3").Value = "'-1.28%"
$ws.Range("G3").Value = "'16"

$ws.Range("D4").Value = "'5.095"
$ws.Range("E4").Value = "'0.86%"
$ws.Range("G4").Value = "'16"

$ws.Range("D5").Value = "'0.08167"
$ws.Range("E5").Value = "'3.70%"
$ws.Range("G5").Value = "'16"

$ws.Range("D6").Value = "'2.091"
$ws.Range("E6").Value = "'-1.45%"
$ws.Range("G6").Value = "'16"

$ws.Range("D7").Value = "'4.145"
$ws.Range("E7").Value = "'0.33%"
$ws.Range("G7").Value = "'16"

$ws.Range("D8").Value = "'7.960"
$ws.Range("E8").Value = "'0.48%"
$ws.Range("G8").Value = "'16"

$ws.Range("D9").Value = "'0.9323"
$ws.Range("E9").Value = "'1.19%"
$ws.Range("G9").Value = "'16"

$ws.Range("D10").Value = "'0.1041"
$ws.Range("E10").Value = "'6.69%"
$ws.Range("G10").Value = "'16"

$ws.Range("E11").Value = "'3.79%"
$ws.Range("G11").Value = "'16"

$ws.Range("D12").Value = "'0.09158"
$ws.Range("E12").Value = "'6.16%"
$ws.Range("G12").Value = "'16"

$ws.Range("D13").Value = "'0.03603"
$ws.Range("E13").Value = "'0.06%"
$ws.Range("G13").Value = "'16"

$ws.Range("D14").Value = "'0.09915"
$ws.Range("E14").Value = "'-0.29%"
$ws.Range("G14").Value = "'16"

$ws.Range("D15").Value = "'0.001433"
$ws.Range("E15").Value = "'-0.54%"
$ws.Range("G15").Value = "'16"

$ws.Range("D16").Value = "'0.005709"
$ws.Range("E16").Value = "'-0.13%"
$ws.Range("G16").Value = "'16"

$ws.Range("D17").Value = "'3.467"
$ws.Range("E17").Value = "'-0.14%"
$ws.Range("G17").Value = "'16"

$ws.Range("D18").Value = "'2.909"
$ws.Range("E18").Value = "'5.66%"
$ws.Range("G18").Value = "'16"

$ws.Range("D19").Value = "'0.3412"
$ws.Range("E19").Value = "'1.12%"
$ws.Range("G19").Value = "'16"

$ws.Range("D20").Value = "'0.1291"
$ws.Range("E20").Value = "'-4.11%"
$ws.Range("G20").Value = "'16"

$ws.Range("D21").Value = "'5.097"
$ws.Range("E21").Value = "'-1.33%"
$ws.Range("G21").Value = "'16"

$ws.Range("G22").Value = "'16"

$ws.Range("D23").Value = "'0.04554"
$ws.Range("E23").Value = "'-0.29%"
$ws.Range("G23").Value = "'16"

$ws.Range("D24").Value = "'0.001243"
$ws.Range("E24").Value = "'0.90%"
$ws.Range("G24").Value = "'16"

$ws.Range("D25").Value = "'0.004798"
$ws.Range("E25").Value = "'-0.43%"
$ws.Range("G25").Value = "'16"

$ws.Range("E26").Value = "'-3.73%"
$ws.Range("G26").Value = "'16"

$ws.Range("D27").Value = "'0.0004502"
$ws.Range("E27").Value = "'-5.24%"
$ws.Range("G27").Value = "'16"

$ws.Range("G28").Value = "'16"

$ws.Range("G29").Value = "'16"

$ws.Range("G30").Value = "'16"

$ws.Range("G31").Value = "'16"

$ws.Range("G32").Value = "'16"

$ws.Range("G33").Value = "'16"

$ws.Range("G34").Value = "'16"

$ws.Range("G35").Value = "'16"

$ws.Range("G36").Value = "'16"

$ws.Range("G37").Value = "'16"

$ws.Range("G38").Value = "'16"

$ws.Range("D39").Value = "'0.01984"
$ws.Range("E39").Value = "'5.64%"
$ws.Range("G39").Value = "'16"

$ws.Range("D40").Value = "'0.04962"
$ws.Range("E40").Value = "'5.44%"
$ws.Range("G40").Value = "'16"

$ws.Range("D41").Value = "'0.007573"
$ws.Range("E41").Value = "'-3.03%"
$ws.Range("G41").Value = "'16"

$ws.Range("D42").Value = "'0.1387"
$ws.Range("E42").Value = "'-0.16%"
$ws.Range("G42").Value = "'16"

$ws.Range("D43").Value = "'0.007872"
$ws.Range("E43").Value = "'1.62%"
$ws.Range("G43").Value = "'16"

$ws.Range("D44").Value = "'0.002222"
$ws.Range("E44").Value = "'0.12%"
$ws.Range("G44").Value = "'16"

$ws.Range("D45").Value = "'0.01153"
$ws.Range("E45").Value = "'1.28%"
$ws.Range("G45").Value = "'16"

$ws.Range("D46").Value = "'0.00006617"
$ws.Range("E46").Value = "'3.65%"
$ws.Range("G46").Value = "'16"

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("G47").Value = "'16"

$ws.Range("D48").Value = "'195.56"
$ws.Range("E48").Value = "'282.57%"
$ws.Range("G48").Value = "'16"

$ws.Range("D49").Value = "'0.001701"
$ws.Range("E49").Value = "'-10.53%"
$ws.Range("G49").Value = "'16"

$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("G50").Value = "'16"

$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.04%"
$ws.Range("G51").Value = "'16"
